# CO_County_Sales_2014_2018.xlsx - fill blank sales/tax cells with "NR"
# on the "aggregate" sheet, and reset the AutoFilter/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aggregate")

# ---------------------------------------------------------------------
# 1. Fill blank cells in columns H:K (for every data row except the
#    "Statewide Total / Not Reported" summary rows) with the text "NR".
#    A reference cell that already carries the desired number style
#    (style used throughout column G/H/I/J/K for normal data rows) is
#    copied onto each blank cell first so the untouched "blank" cells
#    (which still carry an explicit, slightly different empty style)
#    end up matching the style used by cells that never had an explicit
#    style at all.
# ---------------------------------------------------------------------
$styleRef = $ws.Range("G2")
$styleRef.Copy()

$summaryRows = @(66, 131, 196, 261, 326)

for ($r = 2; $r -le 326; $r++) {
    if ($summaryRows -contains $r) { continue }

    foreach ($col in @("H", "I", "J", "K")) {
        $cell = $ws.Range($col + $r)
        $val = $cell.Value()
        if ($val -eq $null) {
            $cell.PasteSpecial(-4122)   # xlPasteFormats
            $cell.Value() = "NR"
        }
    }
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Reset the AutoFilter on "aggregate" so it only spans the header
#    row (A1:K1) instead of the whole data range (A1:K326).
# ---------------------------------------------------------------------
$ws.AutoFilterMode() = $false
$ws.Range("A1:K1").AutoFilter()

$filterName = $wb.Names().Item("aggregate!_FilterDatabase")
$filterName.RefersTo() = "=aggregate!`$A`$1:`$K`$1"

# ---------------------------------------------------------------------
# 3. Update the active selection shown on the "aggregate" sheet.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("O7").Select()
